$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New question/answer rows (29-32) ---------------------------------
# Values are entered in the same order the shared-string table ends up
# in (A29, B29, A30, B30, B31, A31, A32) so the underlying string table
# matches the target edit.
$ws.Range("A29").Value = "How spring framework works?"
$ws.Range("B29").Value = "https://javarevisited.blogspot.com/2017/06/how-spring-mvc-framework-works-web-flow.html"
$ws.Range("A30").Value = "Why do we use ApplicationContext over BeanFactory?"
$ws.Range("B30").Value = "https://dzone.com/articles/difference-between-beanfactory-and-applicationcont"
$ws.Range("B31").Value = "https://howtodoinjava.com/spring-core/spring-bean-post-processors/"
$ws.Range("A31").Value = "What is BeanPostProcessor?"
$ws.Range("A32").Value = "What is BeanFactoryPostProcessor?"

# --- Hyperlinks for the new B-column cells -----------------------------
# Added in row order so the generated relationship ids come out rId12,
# rId13, rId14 (matching the target workbook).
$ws.Hyperlinks.Add($ws.Range("B29"), "https://javarevisited.blogspot.com/2017/06/how-spring-mvc-framework-works-web-flow.html")
$ws.Hyperlinks.Add($ws.Range("B30"), "https://dzone.com/articles/difference-between-beanfactory-and-applicationcont")
$ws.Hyperlinks.Add($ws.Range("B31"), "https://howtodoinjava.com/spring-core/spring-bean-post-processors/")

# --- Re-apply the sheet's normal formatting on top of Hyperlinks.Add's --
# own styling: column A keeps the plain wrap style, column B keeps the
# shared "hyperlink" look already used throughout the sheet (copied from
# an existing hyperlink cell so the same cell style index is reused).
$ws.Range("A28").Copy() | Out-Null
$ws.Range("A29:A32").PasteSpecial(-4122) | Out-Null

$ws.Range("B27").Copy() | Out-Null
$ws.Range("B29:B31").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- Selection / view state ---------------------------------------------
$ws.Range("A33").Select() | Out-Null
